$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D/E columns so numeric-looking strings are not
# auto-converted to numbers by Excel, then restore default style afterwards.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '93.744.74'
$ws.Range("E2").Value = '  -1.62%  '

$ws.Range("D3").Value = '3.323.73'
$ws.Range("E3").Value = '  -3.34%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '230.64'
$ws.Range("E5").Value = '  -3.96%  '

$ws.Range("D6").Value = '617.13'
$ws.Range("E6").Value = '  -3.69%  '

$ws.Range("E7").Value = '  -4.24%  '

$ws.Range("D8").Value = '0.386'
$ws.Range("E8").Value = '  -5.00%  '

$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").Value = '0.929'
$ws.Range("E10").Value = '  -6.18%  '

$ws.Range("D11").Value = '3.322.19'
$ws.Range("E11").Value = '  -3.39%  '

$ws.Range("D12").Value = '41.81'
$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("D13").Value = '0.192'
$ws.Range("E13").Value = '  -3.04%  '

$ws.Range("D14").Value = '93.683.71'
$ws.Range("E14").Value = '  -1.76%  '

$ws.Range("D15").Value = '5.94'
$ws.Range("E15").Value = '  -2.20%  '

$ws.Range("D16").Value = '3.955.07'
$ws.Range("E16").Value = '  -3.52%  '

$ws.Range("D17").Value = '0.0000243'
$ws.Range("E17").Value = '  -4.67%  '

$ws.Range("D18").Value = '8.06'
$ws.Range("E18").Value = '  -4.02%  '

$ws.Range("D19").Value = '3.324.49'
$ws.Range("E19").Value = '  -3.78%  '

$ws.Range("D20").Value = '17.15'
$ws.Range("E20").Value = '  -4.70%  '

$ws.Range("D21").Value = '10.87'
$ws.Range("E21").Value = '  -6.72%  '

$ws.Range("D22").Value = '3.44'
$ws.Range("E22").Value = '  +7.28%  '

$ws.Range("D23").Value = '491.70'
$ws.Range("E23").Value = '  -1.96%  '

$ws.Range("D24").Value = '0.450'
$ws.Range("E24").Value = '  -11.71%  '

$ws.Range("D25").Value = '0.0000181'
$ws.Range("E25").Value = '  -5.34%  '

$ws.Range("D26").Value = '6.03'
$ws.Range("E26").Value = '  -8.45%  '

$ws.Range("D27").Value = '89.51'
$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").Value = '11.64'
$ws.Range("E28").Value = '  -3.60%  '

$ws.Range("D29").Value = '3.510.22'
$ws.Range("E29").Value = '  -3.59%  '

$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("D31").Value = '11.04'
$ws.Range("E31").Value = '  -5.50%  '

$ws.Range("D32").Value = '0.136'
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("D33").Value = '2.62'
$ws.Range("E33").Value = '  -4.23%  '

$ws.Range("D34").Value = '0.990'
$ws.Range("E34").Value = '  -1.15%  '

$ws.Range("D35").Value = '0.173'
$ws.Range("E35").Value = '  -5.32%  '

$ws.Range("D36").Value = '28.21'
$ws.Range("E36").Value = '  -8.12%  '

$ws.Range("D37").Value = '0.528'
$ws.Range("E37").Value = '  -6.98%  '

$ws.Range("D38").Value = '527.23'
$ws.Range("E38").Value = '  +3.54%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").Value = '7.34'
$ws.Range("E40").Value = '  -4.73%  '

$ws.Range("D41").Value = '0.147'
$ws.Range("E41").Value = '  -1.69%  '

$ws.Range("D42").Value = '1.35'
$ws.Range("E42").Value = '  -5.82%  '

$ws.Range("D43").Value = '0.865'
$ws.Range("E43").Value = '  -4.30%  '

$ws.Range("D44").Value = '24.07'
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("E45").Value = '  +5.38%  '

$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = '1.67'
$ws.Range("E46").Value = '  -1.22%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0412'
$ws.Range("E47").Value = '  -0.85%  '

$ws.Range("D48").Value = '5.37'
$ws.Range("E48").Value = '  -2.97%  '

$ws.Range("D49").Value = '52.91'
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("D51").Value = '7.94'
$ws.Range("E51").Value = '  -1.30%  '

# Restore the default (Normal) style on the price/volume columns so no stray
# number-format style is left attached to the cells.
$ws.Range("D2:E51").Style = "Normal"
